$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "DB" contenttype to "Data" (order of first use matches the
# target shared-string table ordering)
$ws.Range("A2").Value = "Data"
$ws.Range("A3").Value = "Data"
$ws.Range("A4").Value = "Data"
$ws.Range("A5").Value = "Data"

# Rename "ML" contenttype to "Machine Learning"
$ws.Range("A6").Value = "Machine Learning"
$ws.Range("A7").Value = "Machine Learning"
$ws.Range("A8").Value = "Machine Learning"

# Extend content descriptions with attribution / data-source info.
# The assignment order below intentionally mirrors the order in which
# the new strings first appear in the saved workbook.
$ws.Range("C6").Value = "Machine Learning Model Predicting Thermoelectric Properties of BiTe-based materials ver.0.2a. Developed by Dr. Jaywan Chung. Engined by LaNN."
$ws.Range("C7").Value = "Machine Learning Model Predicting Thermoelectric Properties of PbTe-based materials ver.0.1a. Developed by Dr. Jaywan Chung. Engined by LaNN."
$ws.Range("C8").Value = "Machine Learning Model Predicting Strength and Conductivity of Copper alloys (POONGSAN data) ver.0.5b. Developed by Dr. Jaywan Chung. Engined by LaNN."
$ws.Range("C10").Value = "Thermoelectric simulator for power TGM, standalone, multistage, 1D. Designed and devloped by Drs. Jaywan Chung and Byungki Ryu."
$ws.Range("C9").Value = "Thermoelectric Power Generation Web Simulator Lite ver.0.53a. Developed by Dr. Jaywan Chung. Designed by Dr. Byungki Ryu."
$ws.Range("C3").Value = "Alloy Design DB (v0.33). Data from KIMS. Designed and developed by Dr. Byungki Ryu"
$ws.Range("C4").Value = "Database of solution, segregation, binding, and doping energies of alloying elements for the design of high-performance Cu-Ni-Si alloys. Data from KIMS. Designed and developed by Dr. Byungki Ryu"
$ws.Range("C5").Value = "A database of calculated solid solution, segregation, and binding characteristics of additive elements in lithium alloys for the design of high-strength, long-life anode materials in Li-S batteries. Data from KIMS. Designed and developed by Dr. Byungki Ryu"
$ws.Range("C2").Value = "Ultra-high quality thermoelectric material property database. Developed by Dr. Byungki Ryu."

# Column D on these rows previously carried a no-op "apply fill" style
# (fillId 0 = no pattern) left over from past editing; clear it so the
# cells fall back to the default style, matching the cleaned-up sheet.
$ws.Range("D2").Interior.Pattern = -4142
$ws.Range("D6").Interior.Pattern = -4142
$ws.Range("D7").Interior.Pattern = -4142
$ws.Range("D8").Interior.Pattern = -4142
$ws.Range("D9").Interior.Pattern = -4142
$ws.Range("D11").Interior.Pattern = -4142

# Update the last active selection to reflect the saved state
$ws.Range("C14").Select()

$wb.Save()
